# Auto-generated edit script applying scheduled price-update diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86,8).Value = 5500  # H86: 5333.3335 -> 5500
$ws.Cells.Item(86,9).Value = 0  # I86: 2000 -> 0
$ws.Cells.Item(86,10).Value = 5500  # J86: 6000 -> 5500
$ws.Cells.Item(86,11).Value = 0  # K86: 2000 -> 0
$ws.Cells.Item(86,12).Value = 5500  # L86: 6000 -> 5500
$ws.Cells.Item(86,13).ClearContents()  # M86: -877 -> (removed)
$ws.Cells.Item(86,14).Value = -7746  # N86: -8246 -> -7746

$ws.Cells.Item(89,8).Value = 5500  # H89: 5333.3335 -> 5500
$ws.Cells.Item(89,9).Value = 0  # I89: 2000 -> 0
$ws.Cells.Item(89,10).Value = 5500  # J89: 6000 -> 5500
$ws.Cells.Item(89,11).Value = 0  # K89: 10000 -> 0
$ws.Cells.Item(89,12).Value = 27500  # L89: 30000 -> 27500
$ws.Cells.Item(89,13).ClearContents()  # M89: -4384 -> (removed)
$ws.Cells.Item(89,14).Value = -38732  # N89: -41232 -> -38732

$ws.Cells.Item(92,8).Value = 1251.4286  # H92: 1023.9286 -> 1251.4286
$ws.Cells.Item(92,9).Value = 655.4545000000001  # I92: 365.9091 -> 655.4545000000001
$ws.Cells.Item(92,11).Value = 655.4545000000001  # K92: 365.9091 -> 655.4545000000001
$ws.Cells.Item(92,13).Value = 592.5454999999999  # M92: 882.0908999999999 -> 592.5454999999999

$ws.Cells.Item(99,8).Value = 1487.2632  # H99: 1430.45 -> 1487.2632
$ws.Cells.Item(99,9).Value = 483.4  # I99: 252.8 -> 483.4
$ws.Cells.Item(99,10).Value = 1845.7858  # J99: 1823 -> 1845.7858
$ws.Cells.Item(99,11).Value = 1450.2  # K99: 758.4000000000001 -> 1450.2
$ws.Cells.Item(99,12).Value = 5537.357400000001  # L99: 5469 -> 5537.357400000001
$ws.Cells.Item(99,13).Value = 47.80000000000018  # M99: 739.5999999999999 -> 47.80000000000018
$ws.Cells.Item(99,14).Value = -8533.357400000001  # N99: -8465 -> -8533.357400000001

$ws.Cells.Item(112,8).Value = 4963.594  # H112: 4966.159 -> 4963.594
$ws.Cells.Item(112,10).Value = 5131.636  # J112: 5134.3184 -> 5131.636
$ws.Cells.Item(112,12).Value = 15394.908  # L112: 15402.9552 -> 15394.908
$ws.Cells.Item(112,14).Value = -17610.908  # N112: -17618.9552 -> -17610.908

$ws.Cells.Item(137,8).Value = 30368378  # H137: 23466954 -> 30368378
$ws.Cells.Item(137,9).Value = 40001930  # I137: 41668636 -> 40001930
$ws.Cells.Item(137,10).Value = 3608516  # J137: 1624932 -> 3608516
$ws.Cells.Item(137,11).Value = 120005790  # K137: 125005908 -> 120005790
$ws.Cells.Item(137,12).Value = 10825548  # L137: 4874796 -> 10825548
$ws.Cells.Item(137,13).Value = -120003240  # M137: -125003358 -> -120003240
$ws.Cells.Item(137,14).Value = -10830648  # N137: -4879896 -> -10830648

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32,8).Value = 29563.486  # H32: 30335.922 -> 29563.486
$ws.Cells.Item(32,9).Value = 29563.486  # I32: 30335.922 -> 29563.486
$ws.Cells.Item(32,11).Value = 29563.486  # K32: 30335.922 -> 29563.486
$ws.Cells.Item(32,13).Value = -29276.486  # M32: -30048.922 -> -29276.486

$ws.Cells.Item(74,8).Value = 2929.35  # H74: 2832.1428 -> 2929.35
$ws.Cells.Item(74,9).Value = 983.61536  # I74: 976.7857 -> 983.61536
$ws.Cells.Item(74,11).Value = 983.61536  # K74: 976.7857 -> 983.61536
$ws.Cells.Item(74,13).Value = -109.61536  # M74: -102.7857 -> -109.61536

$ws.Cells.Item(77,8).Value = 2929.35  # H77: 2832.1428 -> 2929.35
$ws.Cells.Item(77,9).Value = 983.61536  # I77: 976.7857 -> 983.61536
$ws.Cells.Item(77,11).Value = 4918.0768  # K77: 4883.9285 -> 4918.0768
$ws.Cells.Item(77,13).Value = -550.0767999999998  # M77: -515.9285 -> -550.0767999999998

$ws.Cells.Item(102,8).Value = 17035.521  # H102: 20317.736 -> 17035.521
$ws.Cells.Item(102,9).Value = 20149.63  # I102: 23690.812 -> 20149.63
$ws.Cells.Item(102,10).Value = 2243.5  # J102: 2328 -> 2243.5
$ws.Cells.Item(102,11).Value = 20149.63  # K102: 23690.812 -> 20149.63
$ws.Cells.Item(102,12).Value = 2243.5  # L102: 2328 -> 2243.5
$ws.Cells.Item(102,13).Value = -18527.63  # M102: -22068.812 -> -18527.63
$ws.Cells.Item(102,14).Value = -5487.5  # N102: -5572 -> -5487.5

$ws.Cells.Item(110,8).Value = 2116.6  # H110: 2196.3572 -> 2116.6
$ws.Cells.Item(110,9).Value = 1249.8  # I110: 1312.25 -> 1249.8
$ws.Cells.Item(110,11).Value = 1249.8  # K110: 1312.25 -> 1249.8
$ws.Cells.Item(110,13).Value = 795.2  # M110: 732.75 -> 795.2

$ws.Cells.Item(122,8).Value = 3347.4546  # H122: 3570.2 -> 3347.4546
$ws.Cells.Item(122,9).Value = 2643.7334  # I122: 2878.1538 -> 2643.7334
$ws.Cells.Item(122,11).Value = 7931.2002  # K122: 8634.4614 -> 7931.2002
$ws.Cells.Item(122,13).Value = -5481.2002  # M122: -6184.4614 -> -5481.2002

$ws.Cells.Item(132,8).Value = 1252479.1  # H132: 1431191.1 -> 1252479.1
$ws.Cells.Item(132,9).Value = 1431047.6  # I132: 1669306.4 -> 1431047.6
$ws.Cells.Item(132,11).Value = 4293142.800000001  # K132: 5007919.199999999 -> 4293142.800000001
$ws.Cells.Item(132,13).Value = -4290612.800000001  # M132: -5005389.199999999 -> -4290612.800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31,8).Value = 27578.666  # H31: 29051.883 -> 27578.666
$ws.Cells.Item(31,9).Value = 11898.223  # I31: 11838.889 -> 11898.223
$ws.Cells.Item(31,10).Value = 43259.11  # J31: 48416.5 -> 43259.11
$ws.Cells.Item(31,11).Value = 11898.223  # K31: 11838.889 -> 11898.223
$ws.Cells.Item(31,12).Value = 43259.11  # L31: 48416.5 -> 43259.11
$ws.Cells.Item(31,13).Value = -11603.223  # M31: -11543.889 -> -11603.223
$ws.Cells.Item(31,14).Value = -43849.11  # N31: -49006.5 -> -43849.11

$ws.Cells.Item(34,8).Value = 27578.666  # H34: 29051.883 -> 27578.666
$ws.Cells.Item(34,9).Value = 11898.223  # I34: 11838.889 -> 11898.223
$ws.Cells.Item(34,10).Value = 43259.11  # J34: 48416.5 -> 43259.11
$ws.Cells.Item(34,11).Value = 11898.223  # K34: 11838.889 -> 11898.223
$ws.Cells.Item(34,12).Value = 43259.11  # L34: 48416.5 -> 43259.11
$ws.Cells.Item(34,13).Value = -11696.223  # M34: -11636.889 -> -11696.223
$ws.Cells.Item(34,14).Value = -43663.11  # N34: -48820.5 -> -43663.11

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(10,8).Value = 145.75  # H10: 149.16667 -> 145.75
$ws.Cells.Item(10,9).Value = 77.85714  # I10: 85.666664 -> 77.85714
$ws.Cells.Item(10,10).Value = 240.8  # J10: 212.66667 -> 240.8
$ws.Cells.Item(10,11).Value = 233.57142  # K10: 256.999992 -> 233.57142
$ws.Cells.Item(10,12).Value = 722.4000000000001  # L10: 638.00001 -> 722.4000000000001
$ws.Cells.Item(10,13).Value = -94.57141999999999  # M10: -117.999992 -> -94.57141999999999
$ws.Cells.Item(10,14).Value = -1000.4  # N10: -916.00001 -> -1000.4

$ws.Cells.Item(24,8).Value = 637.1429000000001  # H24: 740 -> 637.1429000000001
$ws.Cells.Item(24,10).Value = 715.25  # J24: 947 -> 715.25
$ws.Cells.Item(24,12).Value = 2145.75  # L24: 2841 -> 2145.75
$ws.Cells.Item(24,14).Value = -2605.75  # N24: -3301 -> -2605.75

$ws.Cells.Item(103,8).Value = 2159.1667  # H103: 1000.5714 -> 2159.1667
$ws.Cells.Item(103,9).Value = 1391  # I103: 1000.5714 -> 1391
$ws.Cells.Item(103,10).Value = 6000  # J103: 0 -> 6000
$ws.Cells.Item(103,11).Value = 4173  # K103: 3001.7142 -> 4173
$ws.Cells.Item(103,12).Value = 18000  # L103: 0 -> 18000
$ws.Cells.Item(103,13).Value = -3294  # M103: -2122.7142 -> -3294
$ws.Cells.Item(103,14).Value = -19758  # N103: None -> -19758

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102,8).Value = 2865.9111  # H102: 2908.318 -> 2865.9111
$ws.Cells.Item(102,9).Value = 1939.6428  # I102: 1974.4445 -> 1939.6428
$ws.Cells.Item(102,11).Value = 1939.6428  # K102: 1974.4445 -> 1939.6428
$ws.Cells.Item(102,13).Value = -317.6428000000001  # M102: -352.4445000000001 -> -317.6428000000001

$ws.Cells.Item(126,8).Value = 669908.5600000001  # H126: 727848.4399999999 -> 669908.5600000001
$ws.Cells.Item(126,9).Value = 1113594.4  # I126: 1113647.6 -> 1113594.4
$ws.Cells.Item(126,10).Value = 4379.9  # J126: 4474.875 -> 4379.9
$ws.Cells.Item(126,11).Value = 3340783.2  # K126: 3340942.8 -> 3340783.2
$ws.Cells.Item(126,12).Value = 13139.7  # L126: 13424.625 -> 13139.7
$ws.Cells.Item(126,13).Value = -3338313.2  # M126: -3338472.8 -> -3338313.2
$ws.Cells.Item(126,14).Value = -18079.7  # N126: -18364.625 -> -18079.7

$ws.Cells.Item(135,8).Value = 90538.8  # H135: 85523.28999999999 -> 90538.8
$ws.Cells.Item(135,9).Value = 75899  # I135: 63949.5 -> 75899
$ws.Cells.Item(135,10).Value = 94198.75  # J135: 94152.8 -> 94198.75
$ws.Cells.Item(135,11).Value = 75899  # K135: 63949.5 -> 75899
$ws.Cells.Item(135,12).Value = 94198.75  # L135: 94152.8 -> 94198.75
$ws.Cells.Item(135,13).Value = -70829  # M135: -58879.5 -> -70829
$ws.Cells.Item(135,14).Value = -104338.75  # N135: -104292.8 -> -104338.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(39,8).Value = 10000  # H39: 0 -> 10000
$ws.Cells.Item(39,10).Value = 10000  # J39: 0 -> 10000
$ws.Cells.Item(39,12).Value = 10000  # L39: 0 -> 10000
$ws.Cells.Item(39,14).Value = -10920  # N39: None -> -10920

$ws.Cells.Item(40,8).Value = 4380.8  # H40: 3922.5557 -> 4380.8
$ws.Cells.Item(40,9).Value = 4380.8  # I40: 3922.5557 -> 4380.8
$ws.Cells.Item(40,11).Value = 4380.8  # K40: 3922.5557 -> 4380.8
$ws.Cells.Item(40,13).Value = -4244.8  # M40: -3786.5557 -> -4244.8

$ws.Cells.Item(93,8).Value = 1428.0454  # H93: 1560.1818 -> 1428.0454
$ws.Cells.Item(93,9).Value = 1254.1  # I93: 1544.7 -> 1254.1
$ws.Cells.Item(93,10).Value = 1573  # J93: 1573.0834 -> 1573
$ws.Cells.Item(93,11).Value = 1254.1  # K93: 1544.7 -> 1254.1
$ws.Cells.Item(93,12).Value = 1573  # L93: 1573.0834 -> 1573
$ws.Cells.Item(93,13).Value = -6.099999999999909  # M93: -296.7 -> -6.099999999999909
$ws.Cells.Item(93,14).Value = -4069  # N93: -4069.0834 -> -4069

$ws.Cells.Item(122,8).Value = 3741.0286  # H122: 3758.1714 -> 3741.0286
$ws.Cells.Item(122,9).Value = 3582.4827  # I122: 3603.1724 -> 3582.4827
$ws.Cells.Item(122,11).Value = 10747.4481  # K122: 10809.5172 -> 10747.4481
$ws.Cells.Item(122,13).Value = -8297.4481  # M122: -8359.5172 -> -8297.4481

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(64,8).Value = 45000  # H64: 0 -> 45000
$ws.Cells.Item(64,10).Value = 45000  # J64: 0 -> 45000
$ws.Cells.Item(64,12).Value = 45000  # L64: 0 -> 45000
$ws.Cells.Item(64,14).Value = -45496  # N64: None -> -45496

$ws.Cells.Item(67,8).Value = 45000  # H67: 0 -> 45000
$ws.Cells.Item(67,10).Value = 45000  # J67: 0 -> 45000
$ws.Cells.Item(67,12).Value = 45000  # L67: 0 -> 45000
$ws.Cells.Item(67,14).Value = -46716  # N67: None -> -46716

$ws.Cells.Item(107,8).Value = 1576.6471  # H107: 1659.5625 -> 1576.6471
$ws.Cells.Item(107,9).Value = 1017.75  # I107: 1087.5454 -> 1017.75
$ws.Cells.Item(107,11).Value = 3053.25  # K107: 3262.6362 -> 3053.25
$ws.Cells.Item(107,13).Value = -1133.25  # M107: -1342.6362 -> -1133.25

$ws.Cells.Item(113,8).Value = 2919.875  # H113: 3111.2 -> 2919.875
$ws.Cells.Item(113,9).Value = 550.5714  # I113: 634 -> 550.5714
$ws.Cells.Item(113,11).Value = 1651.7142  # K113: 1902 -> 1651.7142
$ws.Cells.Item(113,13).Value = 518.2857999999999  # M113: 268 -> 518.2857999999999

$ws.Cells.Item(122,8).Value = 1826.4193  # H122: 1860.8 -> 1826.4193
$ws.Cells.Item(122,10).Value = 2436.125  # J122: 2670.5715 -> 2436.125
$ws.Cells.Item(122,12).Value = 7308.375  # L122: 8011.7145 -> 7308.375
$ws.Cells.Item(122,14).Value = -12208.375  # N122: -12911.7145 -> -12208.375

$ws.Cells.Item(132,8).Value = 23605582  # H132: 26752448 -> 23605582
$ws.Cells.Item(132,9).Value = 25080430  # I132: 28662764 -> 25080430
$ws.Cells.Item(132,11).Value = 75241290  # K132: 85988292 -> 75241290
$ws.Cells.Item(132,13).Value = -75238760  # M132: -85985762 -> -75238760

$ws.Cells.Item(136,8).Value = 18531800  # H136: 20847702 -> 18531800
$ws.Cells.Item(136,9).Value = 20847088  # I136: 22236748 -> 20847088
$ws.Cells.Item(136,10).Value = 9484.5  # J136: 12000 -> 9484.5
$ws.Cells.Item(136,11).Value = 62541264  # K136: 66710244 -> 62541264
$ws.Cells.Item(136,12).Value = 28453.5  # L136: 36000 -> 28453.5
$ws.Cells.Item(136,13).Value = -62538714  # M136: -66707694 -> -62538714
$ws.Cells.Item(136,14).Value = -33553.5  # N136: -41100 -> -33553.5
